# Apply crypto price/volume updates per commit "Updated cryptos list"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '51.879.14'
$ws.Range('E2').Value = '  -0.38%  '
$ws.Range('D3').Value = '2.820.82'
$ws.Range('E3').Value = '  +1.29%  '
$ws.Range('D4').Value = "'1.00"
$ws.Range('E4').Value = '  +0.11%  '
$ws.Range('D5').Value = "'355.64"
$ws.Range('E5').Value = '  +3.79%  '
$ws.Range('D6').Value = "'112.13"
$ws.Range('E6').Value = '  -2.98%  '
$ws.Range('D7').Value = "'0.564"
$ws.Range('E7').Value = '  +2.92%  '
$ws.Range('E8').Value = '  +0.06%  '
$ws.Range('D9').Value = "'0.601"
$ws.Range('E9').Value = '  +3.94%  '
$ws.Range('D10').Value = "'40.95"
$ws.Range('E10').Value = '  -2.78%  '
$ws.Range('D11').Value = "'0.0856"
$ws.Range('E11').Value = '  -0.15%  '
$ws.Range('B12').Value = 'TRON'
$ws.Range('C12').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D12').Value = "'0.132"
$ws.Range('E12').Value = '  +1.26%  '
$ws.Range('B13').Value = 'Chainlink'
$ws.Range('C13').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D13').Value = "'19.97"
$ws.Range('E13').Value = '  -0.17%  '
$ws.Range('D14').Value = "'7.75"
$ws.Range('E14').Value = '  +1.55%  '
$ws.Range('D15').Value = '3.270.53'
$ws.Range('E15').Value = '  +1.65%  '
$ws.Range('D16').Value = '2.838.74'
$ws.Range('E16').Value = '  +2.26%  '
$ws.Range('D17').Value = "'0.926"
$ws.Range('E17').Value = '  +5.48%  '
$ws.Range('D18').Value = '51.860.77'
$ws.Range('E18').Value = '  -0.15%  '
$ws.Range('D20').Value = "'3.14"
$ws.Range('E20').Value = '  -1.93%  '
$ws.Range('D21').Value = "'13.40"
$ws.Range('E21').Value = '  +1.30%  '
$ws.Range('D22').Value = '0.0₃0992'
$ws.Range('E22').Value = '  +1.22%  '
$ws.Range('D23').Value = "'69.91"
$ws.Range('E23').Value = '  -0.31%  '
$ws.Range('D24').Value = "'267.98"
$ws.Range('E24').Value = '  -3.30%  '
$ws.Range('E25').Value = '  +1.38%  '
$ws.Range('D26').Value = "'27.02"
$ws.Range('E26').Value = '  +1.24%  '
$ws.Range('E27').Value = '  +0.09%  '
$ws.Range('D28').Value = "'10.31"
$ws.Range('E28').Value = '  +1.06%  '
$ws.Range('E29').Value = '  +1.48%  '
$ws.Range('B30').Value = 'VeChain'
$ws.Range('C30').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D30').Value = "'0.0476"
$ws.Range('E30').Value = '  +24.57%  '
$ws.Range('B31').Value = 'Kaspa'
$ws.Range('C31').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D31').Value = "'0.140"
$ws.Range('E31').Value = '  -1.00%  '
$ws.Range('B32').Value = 'OKB'
$ws.Range('C32').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D32').Value = "'52.57"
$ws.Range('E32').Value = '  +4.62%  '
$ws.Range('B33').Value = 'InjectiveProtocol'
$ws.Range('C33').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D33').Value = "'34.93"
$ws.Range('E33').Value = '  +0.31%  '
$ws.Range('D34').Value = "'5.89"
$ws.Range('E34').Value = '  +3.07%  '
$ws.Range('D35').Value = "'5.39"
$ws.Range('E35').Value = '  +8.44%  '
$ws.Range('D36').Value = "'0.0843"
$ws.Range('E36').Value = '  +2.93%  '
$ws.Range('E37').Value = '  -0.08%  '
$ws.Range('E38').Value = '  +1.94%  '
$ws.Range('B39').Value = 'ARBITRUM'
$ws.Range('C39').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D39').Value = "'2.05"
$ws.Range('E39').Value = '  -2.88%  '
$ws.Range('B40').Value = 'Celestia'
$ws.Range('C40').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D40').Value = "'18.48"
$ws.Range('E40').Value = '  -2.63%  '
$ws.Range('D41').Value = "'0.116"
$ws.Range('E41').Value = '  +0.35%  '
$ws.Range('E42').Value = '  -4.79%  '
$ws.Range('D43').Value = "'23.28"
$ws.Range('E43').Value = '  -0.42%  '
$ws.Range('D44').Value = "'124.46"
$ws.Range('E44').Value = '  -1.69%  '
$ws.Range('D45').Value = "'2.27"
$ws.Range('E45').Value = '  -2.81%  '
$ws.Range('B46').Value = 'Maker'
$ws.Range('C46').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D46').Value = '2.095.20'
$ws.Range('E46').Value = '  +1.07%  '
$ws.Range('B47').Value = 'NEARProtocol'
$ws.Range('C47').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D47').Value = "'3.35"
$ws.Range('E47').Value = '  +0.80%  '
$ws.Range('D49').Value = "'5.98"
$ws.Range('E49').Value = '  +7.46%  '
$ws.Range('D50').Value = "'0.971"
$ws.Range('E50').Value = '  +7.60%  '
$ws.Range('E51').Value = '  +2.12%  '
